$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet to reflect new export timestamp
$ws.Name = "IClientBalance-20240703-104417-"

# Update reference date (column G) for all data rows from 2024-06-28 (45471) to 2024-07-03 (45476)
$dateRange = $ws.Range("G2:G275")
$dateRange.Value = 45476

# Row-specific balance corrections (Saldo Previsto / Vl. Total, and one Vl. Projetado)
$ws.Range("E5").Value = 203.1
$ws.Range("H5").Value = 203.1
$ws.Range("E8").Value = 522.07
$ws.Range("H8").Value = 522.07
$ws.Range("E43").Value = 1864.5
$ws.Range("H43").Value = 1864.5
$ws.Range("E52").Value = 178319.92
$ws.Range("H52").Value = 178319.92
$ws.Range("E55").Value = 128119.55
$ws.Range("H55").Value = 128119.55
$ws.Range("E57").Value = 63513.86
$ws.Range("H57").Value = 63513.86
$ws.Range("E60").Value = 892.69
$ws.Range("H60").Value = 892.69
$ws.Range("E99").Value = 119
$ws.Range("H99").Value = 119
$ws.Range("E104").Value = 426.4
$ws.Range("H104").Value = 426.4
$ws.Range("E108").Value = 261.77
$ws.Range("H108").Value = 261.77
$ws.Range("E109").Value = 7398.85
$ws.Range("H109").Value = 7398.85
$ws.Range("E110").Value = 3521.32
$ws.Range("H110").Value = 3521.32
$ws.Range("E112").Value = -51.27
$ws.Range("H112").Value = -51.27
$ws.Range("E114").Value = 187.48
$ws.Range("H114").Value = 187.48
$ws.Range("E118").Value = 2802.3
$ws.Range("H118").Value = 2802.3
$ws.Range("E138").Value = 86.35
$ws.Range("H138").Value = 86.35
$ws.Range("E143").Value = 912.41
$ws.Range("H143").Value = 912.41
$ws.Range("E148").Value = 21457.78
$ws.Range("H148").Value = 21457.78
$ws.Range("D158").Value = 0
$ws.Range("E158").Value = 6907.1
$ws.Range("H158").Value = 6907.1
$ws.Range("E161").Value = 246.23
$ws.Range("H161").Value = 246.23
$ws.Range("E165").Value = 61788.61
$ws.Range("H165").Value = 61788.61
$ws.Range("E171").Value = 0
$ws.Range("H171").Value = 0
$ws.Range("E172").Value = 10933.97
$ws.Range("H172").Value = 10933.97
$ws.Range("E213").Value = 1072.56
$ws.Range("H213").Value = 1072.56
$ws.Range("E224").Value = 435.85
$ws.Range("H224").Value = 435.85
$ws.Range("E230").Value = 33404.28
$ws.Range("H230").Value = 33404.28
$ws.Range("E232").Value = 5.9
$ws.Range("H232").Value = 5.9
$ws.Range("E235").Value = 95.58
$ws.Range("H235").Value = 95.58
$ws.Range("E249").Value = 101.53
$ws.Range("H249").Value = 101.53
$ws.Range("E255").Value = 47458.32
$ws.Range("H255").Value = 47458.32
$ws.Range("E264").Value = 817.94
$ws.Range("H264").Value = 817.94
$ws.Range("E270").Value = 1237.52
$ws.Range("H270").Value = 1237.52
$ws.Range("E271").Value = 274.26
$ws.Range("H271").Value = 274.26
$ws.Range("E272").Value = 194.57
$ws.Range("H272").Value = 194.57
$ws.Range("E274").Value = 412.58
$ws.Range("H274").Value = 412.58
